$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C (rows 2-28) holds a "Förändrad" (changed) date that needs to
# move from 45478 (2024-07-05) to 45479 (2024-07-06).
for ($row = 2; $row -le 28; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $cell.Value = 45479
}
